# Update the March-2025 attendance sheet:
#  - rename the "Entry Time" / "Exit Time" headers to "EntryTime" / "ExitTime"
#  - append the new attendance record for Vasanth Kumar on row 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "EntryTime"
$ws.Range("D1").Value = "ExitTime"

# Force the date-looking value to be stored as plain text (matches the
# source data export), then drop back to the default "Normal" style so
# no stray number-format style is left behind on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-03-09"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "Vasanth Kumar"
$ws.Range("C2").Value = "15:49:49"
$ws.Range("D2").Value = "15:51:05"
